$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N10").Value = 9684.19
$ws.Range("O10").Value = 9404.28

$ws.Range("N12").Value = 493468.88

$ws.Range("M13").Value = 524661.53
$ws.Range("N13").Value = 168003.02
$ws.Range("O13").Value = 148022.19

$ws.Range("N15").Value = 645.18

$ws.Range("K17").Value = 10198.44

$ws.Range("N18").Value = 4732.5
$ws.Range("O18").Value = 4732.5

$ws.Range("N25").Value = 105087.66
$ws.Range("O25").Value = 105087.66
